$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Countries table: swap Estonia / Azerbaiyan rows (labels only; the
# underlying per-row stats below are updated in place to match) ---
$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("A74").Value = "Estonia"

# --- Row 41 (country index 45) updated stats ---
$ws.Range("D41").Value = 5959
$ws.Range("E41").Value = 2312
$ws.Range("F41").Value = 72
$ws.Range("G41").Value = 5
$ws.Range("H41").Value = 427

# --- Row 42 (country index 46) updated stats ---
$ws.Range("B42").Value = 8275
$ws.Range("C42").Value = 233
$ws.Range("D42").Value = 1209
$ws.Range("E42").Value = 6904
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 162

# --- Row 73 (now Azerbaiyan) updated stats ---
$ws.Range("B73").Value = 1678
$ws.Range("C73").Value = 33
$ws.Range("D73").Value = 1162
$ws.Range("E73").Value = 494
$ws.Range("F73").Value = 15
$ws.Range("H73").Value = 22

# --- Row 74 (now Estonia) updated stats ---
$ws.Range("B74").Value = 1647
$ws.Range("C74").Value = 4
$ws.Range("D74").Value = 233
$ws.Range("E74").Value = 1364
$ws.Range("F74").Value = 7
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 50

# --- Row 83 updated stat ---
$ws.Range("F83").Value = 12

# --- Row 105 updated stats ---
$ws.Range("B105").Value = 571
$ws.Range("C105").Value = 48
$ws.Range("E105").Value = 438
